# Update the timestamp portion of the test-data email addresses from
# 20251109_003734 to 20251109_004215. The same underlying string values
# are reused (shared) by both the "UsuariosRegistro" sheet (column C,
# rows 2-6) and the "LoginData" sheet (column A, rows 2-3), so both need
# to be kept in sync.

$oldStamp = "20251109_003734"
$newStamp = "20251109_004215"

$wb = $excel.ActiveWorkbook

$wsRegistro = $wb.Worksheets.Item("UsuariosRegistro")
for ($row = 2; $row -le 6; $row++) {
    $cell = $wsRegistro.Cells.Item($row, 3)
    $current = $cell.Value()
    $updated = $current.Replace($oldStamp, $newStamp)
    $cell.Value = $updated
}

$wsLogin = $wb.Worksheets.Item("LoginData")
for ($row = 2; $row -le 3; $row++) {
    $cell = $wsLogin.Cells.Item($row, 1)
    $current = $cell.Value()
    $updated = $current.Replace($oldStamp, $newStamp)
    $cell.Value = $updated
}
